# fix shark double counts
# Update "Coverage (%) Update" (column C) values on Sheet1 for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("C2").Value  = 88.70339407047719
$ws.Range("C9").Value  = 95.38283744105355
$ws.Range("C10").Value = 95.35980609945506
$ws.Range("C11").Value = 98.4871602932466
$ws.Range("C14").Value = 98.26681293112539
$ws.Range("C15").Value = 87.65426515699606
